$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.907933235168457
$ws.Range("B1").Value = 2.778306722640991
$ws.Range("C1").Value = 3.410005807876587
$ws.Range("D1").Value = 1.116803884506226
$ws.Range("E1").Value = 0.7195398807525635
